$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values to re-pulled/recalculated data
$ws.Range("F2").Value = -8
$ws.Range("F16").Value = 3
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = 3
$ws.Range("F34").Value = -2
$ws.Range("F35").Value = 0
$ws.Range("F36").Value = 3
$ws.Range("F38").Value = -1
$ws.Range("F39").Value = 3
$ws.Range("F46").Value = -3
$ws.Range("F53").Value = 0
